$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the front (A, B). Everything that was in
# A:P shifts right to C:R.
$ws.Range("A1:B1").EntireColumn.Insert()

# New header labels for the two inserted columns.
$ws.Range("A1").Value = "point_number_OPX"
$ws.Range("B1").Value = "point_number_CPX"

# New per-row "point number" labels for rows 2-16.
$pointNumOPX = @("2 / 1 . ", "2 / 2 . ", "2 / 3 . ", "2 / 4 . ", "2 / 5 . ", "2 / 6 . ", "2 / 7 . ", "2 / 8 . ", "2 / 9 . ", "2 / 10 . ", "2 / 11 . ", "2 / 12 . ", "2 / 13 . ", "2 / 14 . ", "2 / 15 . ")
$pointNumCPX = @("3 / 1 . ", "3 / 2 . ", "3 / 3 . ", "3 / 4 . ", "3 / 5 . ", "3 / 6 . ", "3 / 7 . ", "3 / 8 . ", "3 / 9 . ", "3 / 10 . ", "3 / 11 . ", "3 / 12 . ", "3 / 13 . ", "3 / 14 . ", "3 / 15 . ")

for ($i = 0; $i -lt $pointNumOPX.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $pointNumOPX[$i]
    $ws.Cells.Item($row, 2).Value = $pointNumCPX[$i]
}

# Center the new data cells, matching the look of the other "computed"
# columns (e.g. the Nb_ions_Ca_CPX column) elsewhere in the sheet.
$dataRange = $ws.Range("A2:B16")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# Resize the columns to fit the new (wider) header/data text - mirrors
# the widths left behind by Excel's auto-fit after the edit.
$ws.Columns.Item(1).ColumnWidth = 18.592447916666668
$ws.Columns.Item(2).ColumnWidth = 17.736979166666668
$ws.Columns.Item(3).ColumnWidth = 17.307291666666668

# The selection left behind after the edit.
$ws.Range("D31").Select()
